# Update "Descrição das Fontes de Recursos" worksheet with revised
# classification entries (Rio Doce repactuação, FUNDEB VAAR, and the
# split of "Transferências Especiais de Recursos da União" into a
# pre-2023 and a post-2024 variant).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 81 - was "CONTRIBUIÇÃO DO SERVIDOR PARA O FUNPREV"
$ws.Range("B81").Value = "RECURSOS DO ACORDO DE REPACTUAÇÃO DO RIO DOCE"
$ws.Range("C81").Value = "RECURSOS PROVENIENTES DO ACORDO JUDICIAL DE REPARAÇÃO INTEGRAL E DEFINITIVA RELATIVA AO ROMPIMENTO DA BARRAGEM DE FUNDÃO EM MARIANA."

# Row 91 - was "RECURSOS DECORRENTES DA COBRANÇA DOS CRÉDITOS INADIMPLIDOS INSCRITOS - LEI 22.606/2017"
$ws.Range("B91").Value = "TRANSFERÊNCIAS DO FUNDEB - COMPLEMENTAÇÃO DA UNIÃO - VAAR"
$ws.Range("C91").Value = "CONTROLE DOS RECURSOS DE COMPLEMENTAÇÃO DA UNIÃO AO FUNDEB - VAAR, COM BASE NA ALÍNEA C, INCISO V DO ART. 212-A DA CONSTITUIÇÃO FEDERAL."

# Row 98 - was "TRANSFERÊNCIAS ESPECIAIS DE RECURSOS DA UNIÃO" (now the "até 2023" variant)
$ws.Range("B98").Value = "TRANSFERÊNCIAS ESPECIAIS DE RECURSOS DA UNIÃO - RECURSOS RECEBIDOS ATÉ 2023"
$ws.Range("C98").Value = "RECURSOS TRANSFERIDOS PELA UNIÃO NA FORMA ESTABELECIDA NO INCISO I DO ART.166-A DA CONSTITUIÇÃO DA REPÚBLICA FEDERATIVA DO BRASIL DE 1988. REFERE-SE AOS RECURSOS TRANSFERIDOS PELA UNIÃO AO ESTADO ATÉ O FINAL DO EXERCÍCIO FINANCEIRO DE 2023, BEM COMO A SUAS RESPECTIVAS RESTITUIÇÕES E REMUNERAÇÕES DE DEPÓSITOS BANCÁRIOS, AINDA QUE REGISTRADAS EM EXERCÍCIOS FINANCEIROS POSTERIORES."

# Row 99 - was "CONTRIBUIÇÃO MILITAR PARA CUSTEIO DO BENEFÍCIO DE ASSISTÊNCIA À SAÚDE DOS MILITARES E DEPENDENTES"
# (now re-purposed as the new "a partir de 2024" variant of Transferências Especiais)
$ws.Range("B99").Value = "TRANSFERÊNCIAS ESPECIAIS DE RECURSOS DA UNIÃO"
$ws.Range("C99").Value = "RECURSOS TRANSFERIDOS PELA UNIÃO NA FORMA ESTABELECIDA NO INCISO I DO ART.166-A DA CONSTITUIÇÃO DA REPÚBLICA FEDERATIVA DO BRASIL DE 1988. REFERE-SE AOS RECURSOS TRANSFERIDOS PELA UNIÃO AO ESTADO A PARTIR DO EXERCÍCIO FINANCEIRO DE 2024, BEM COMO A SUAS RESPECTIVAS RESTITUIÇÕES E REMUNERAÇÕES DE DEPÓSITOS BANCÁRIOS."

# Reflect the author's final selection state (cell C81, the last entry edited)
$ws.Range("C81").Select()
